# Apply updated odds values to Sheet1 as described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 3.2
$ws.Range("I2").Value = 2.5
$ws.Range("J2").Value = 3.75
$ws.Range("L2").Value = 3.2
$ws.Range("X2").Value = 15
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 34
$ws.Range("AA2").Value = 29
$ws.Range("AH2").Value = 7
$ws.Range("AI2").Value = 11
$ws.Range("AK2").Value = 23
$ws.Range("AL2").Value = 21

# Row 4 updates
$ws.Range("Q4").Value = 2.35
$ws.Range("R4").Value = 1.57

# Row 6 updates
$ws.Range("H6").Value = 3.2
$ws.Range("J6").Value = 3.25
$ws.Range("N6").Value = 7.5
$ws.Range("U6").Value = 1.91
$ws.Range("V6").Value = 1.8
$ws.Range("AA6").Value = 23
$ws.Range("AH6").Value = 7.5
$ws.Range("AJ6").Value = 11
$ws.Range("AK6").Value = 29
$ws.Range("AL6").Value = 26
$ws.Range("AN6").Value = 4.5
$ws.Range("AS6").Value = 201
$ws.Range("AX6").Value = 17
